# Generate Report for Archive
#
# 1. Update status text "Ready for handoff" -> "In Translation" on the
#    Overview sheet (columns mirroring zh-cn/de-de status) and on the
#    per-locale "Status" column of the zh-cn / de-de sheets.
# 2. Narrow the "Status" column (and the Overview sheet's mirrored zh-cn /
#    de-de columns) so they are no longer as wide.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Update the status text wherever it currently reads "Ready for handoff" ---
$statusCells = @($overview.Range("E2"), $overview.Range("F2"), $zhcn.Range("C2"), $dede.Range("C2"))
foreach ($cell in $statusCells) {
    if ($cell.Value2 -eq "Ready for handoff") {
        $cell.Value2 = "In Translation"
    }
}

# --- Resize columns (Status / zh-cn / de-de) ---
# Target character width is ~13.41; Excel's ColumnWidth setter snaps to the
# screen-pixel grid of the workbook's Normal-style font, so 12.5 is the
# input that lands on the closest reachable grid value.
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
